# Update "Detailed technical stack.docx":
#   - Turn the "GO into details" paragraph into a bold Heading 1 paragraph
#     reading "Front-end (Drag&Drop UI) detailed development with Angular<tab>"
#   - Normalize the trailing empty paragraph so it carries an explicit
#     en-US language run property.

$d = $word.ActiveDocument

# Locate the target paragraph ("GO into details") by its text content,
# and the blank paragraph that immediately follows it, rather than relying
# on a hard-coded paragraph index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    $ptext = $ptext.TrimEnd([char]13, [char]7)
    if ($ptext -eq "GO into details") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'GO into details' paragraph"
}

$pHeading = $d.Paragraphs.Item($targetIndex)
$pBlank = $d.Paragraphs.Item($targetIndex + 1)

$rangeStart = $pHeading.Range.Start
$rangeEnd = $pBlank.Range.End

$target = $d.Range($rangeStart, $rangeEnd)

# New heading paragraph: bold Heading 1 style, split into the same runs
# (and spell-check markers around "Drag&Drop") as produced by a live edit,
# followed by the normalized blank paragraph.
$headingParagraphXml = '<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Front</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>end (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Drag&amp;Drop</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> UI) </w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">detailed </w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>development wi</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">th </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Angular</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r>' + `
    '</w:p>'

$blankParagraphXml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'

$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $headingParagraphXml + $blankParagraphXml + '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$target.InsertXML($payload)
